$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing remark cell (C2) with extra text
$ws.Range("C2").Value = "不需要语言输入，直接输入视频，且叙述者不需要标注每个活动。使用了GMM"

# Add row 3: Distance Based Ranking Models
$ws.Range("A3").Value = "Distance Based Ranking Models.  "
$ws.Range("B3").Value = "M. A. Fligner and J. S. Verducci"

# Add row 4: Action Recognition by Dense Trajectories
$ws.Range("A4").Value = "Action Recognition by Dense Trajectories"
$ws.Range("B4").Value = "Heng Wang, Alexander Kläser, Cordelia Schmid, Liu Cheng-Lin`n"
$ws.Range("C4").Value = "密集轨迹法，估测物体运动的经典方法。"

# Wrap text on B4 and set the row height to fit the wrapped 2-line text
$ws.Range("B4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 42

# Adjust column B width (closest achievable value to the target 36.58203125 chars
# given this engine's internal pixel-snapping of the ColumnWidth property)
$ws.Columns.Item(2).ColumnWidth = 35.857142857142854

# Update view: top-left cell and selection
$ws.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
